$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "Save"
$ws.Range("H1").Style = $ws.Range("G1").Style

$ws.Range("H2").Value = 0
